$wb = $excel.ActiveWorkbook

# The workbook has 4 sheets: 展览, 演出, 本地生活, 全部类型
# The "展览" and "全部类型" sheets share identical data and both need
# their "想去人数" (F column) counts bumped up.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 115
    $ws.Range("F11").Value = 4453
    $ws.Range("F14").Value = 1239
    $ws.Range("F17").Value = 814
}
